$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 24 de Septiembre de 2020 a las 08:40"

# Israel (row 27)
$ws.Range("A27").Value = "Israel"
$ws.Range("B27").Value = 206332
$ws.Range("C27").Value = 1642
$ws.Range("D27").Value = 148075
$ws.Range("E27").Value = 56922
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = 1335

# Ucrania (row 28)
$ws.Range("A28").Value = "Ucrania"
$ws.Range("B28").Value = 188106
$ws.Range("C28").Value = 3372
$ws.Range("D28").Value = 83458
$ws.Range("E28").Value = 100891
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 52
$ws.Range("H28").Value = 3757

# Uzbekistan (row 60)
$ws.Range("A60").Value = "Uzbekistan"
$ws.Range("B60").Value = 53667
$ws.Range("C60").Value = 392
$ws.Range("D60").Value = 49832
$ws.Range("E60").Value = 3390
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 445

# Kirguistan (row 66)
$ws.Range("A66").Value = "Kirguistan"
$ws.Range("B66").Value = 45757
$ws.Range("C66").Value = 127
$ws.Range("D66").Value = 42005
$ws.Range("E66").Value = 2689
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 1063

# Row 82 now holds Hungria's updated data (Hungria overtakes Camerun)
$ws.Range("A82").Value = "Hungria"
$ws.Range("B82").Value = 21200
$ws.Range("C82").Value = 750
$ws.Range("D82").Value = 4818
$ws.Range("E82").Value = 15673
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 7
$ws.Range("H82").Value = 709

# Row 83 now holds Camerun's (unchanged) data
$ws.Range("A83").Value = "Camerun"
$ws.Range("B83").Value = 20690
$ws.Range("C83").Value = 0
$ws.Range("D83").Value = 19124
$ws.Range("E83").Value = 1150
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 416

# Row 214 now holds Montserrat's data
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

# Row 215 now holds Islas Malvinas' data
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0
